# Add the team's season record (Wins / Losses / Ties) as three new
# columns appended to the right of the existing roster/stat table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column past the previous last column (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used by the rest of the header row
# (bold, bordered, centered) by copying an existing header cell's style.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Season record for the 2009 Colorado Rockies: 92 wins, 70 losses, 0 ties.
# Stamp it on every player row (2 through 47).
$record = @{ "AD" = 92; "AE" = 70; "AF" = 0 }
foreach ($col in $record.Keys) {
    $ws.Range("$($col)2`:$($col)47").Value = $record[$col]
}
